$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 97, shifting existing rows 97:134 down to 98:135
$ws.Rows.Item(97).Insert()

# Populate the newly inserted row 97 with the new data record
$ws.Range("A97").Value = 6
$ws.Range("B97").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 44559
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = "Berries"
$ws.Range("I97").Value = 100101004
$ws.Range("J97").Value = "Frambuesa"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 300
$ws.Range("N97").Value = 7000
$ws.Range("O97").Value = 8000
$ws.Range("P97").Value = 7500
$ws.Range("Q97").Value = "$/bandeja 2 kilos"
$ws.Range("R97").Value = "Provincia de Linares"
$ws.Range("S97").Value = 3750
$ws.Range("T97").Value = 2
